$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (row 386-464), continuing the series through 2021-12-08
$newData = @(
    @(44460,0,8,53.02578378736661),
    @(44461,0,7,46.39756081394578),
    @(44462,1,6,39.76933784052495),
    @(44463,1,7,46.39756081394578),
    @(44464,1,7,46.39756081394578),
    @(44465,2,8,53.02578378736661),
    @(44466,1,6,39.76933784052495),
    @(44467,2,8,53.02578378736661),
    @(44468,0,8,53.02578378736661),
    @(44469,1,8,53.02578378736661),
    @(44470,2,9,59.65400676078744),
    @(44471,3,11,72.91045270762909),
    @(44472,0,9,59.65400676078744),
    @(44473,0,8,53.02578378736661),
    @(44474,2,8,53.02578378736661),
    @(44475,0,8,53.02578378736661),
    @(44476,1,8,53.02578378736661),
    @(44477,0,6,39.76933784052495),
    @(44478,0,3,19.88466892026248),
    @(44479,1,4,26.5128918936833),
    @(44480,0,4,26.5128918936833),
    @(44481,0,2,13.25644594684165),
    @(44482,0,2,13.25644594684165),
    @(44483,0,1,6.628222973420826),
    @(44484,2,3,19.88466892026248),
    @(44485,3,6,39.76933784052495),
    @(44486,0,5,33.14111486710413),
    @(44487,0,5,33.14111486710413),
    @(44488,1,6,39.76933784052495),
    @(44489,0,6,39.76933784052495),
    @(44490,0,6,39.76933784052495),
    @(44491,1,5,33.14111486710413),
    @(44492,1,3,19.88466892026248),
    @(44493,1,4,26.5128918936833),
    @(44494,0,4,26.5128918936833),
    @(44495,0,3,19.88466892026248),
    @(44496,0,3,19.88466892026248),
    @(44497,0,3,19.88466892026248),
    @(44498,0,2,13.25644594684165),
    @(44499,1,2,13.25644594684165),
    @(44500,2,3,19.88466892026248),
    @(44501,1,4,26.5128918936833),
    @(44502,0,4,26.5128918936833),
    @(44503,0,4,26.5128918936833),
    @(44504,0,4,26.5128918936833),
    @(44505,1,5,33.14111486710413),
    @(44506,1,5,33.14111486710413),
    @(44507,0,3,19.88466892026248),
    @(44508,1,3,19.88466892026248),
    @(44509,0,3,19.88466892026248),
    @(44510,0,3,19.88466892026248),
    @(44511,4,7,46.39756081394578),
    @(44512,5,11,72.91045270762909),
    @(44513,2,12,79.53867568104991),
    @(44514,5,17,112.679790548154),
    @(44515,4,20,132.5644594684165),
    @(44516,23,43,285.0135878570956),
    @(44517,0,43,285.0135878570956),
    @(44518,0,39,258.5006959634122),
    @(44519,1,35,231.9878040697289),
    @(44520,1,34,225.3595810963081),
    @(44521,4,33,218.7313581228873),
    @(44522,4,33,218.7313581228873),
    @(44523,2,12,79.53867568104991),
    @(44524,10,22,145.8209054152582),
    @(44525,4,26,172.3337973089415),
    @(44526,3,28,185.5902432557831),
    @(44527,1,28,185.5902432557831),
    @(44528,6,30,198.8466892026248),
    @(44529,1,27,178.9620202823623),
    @(44530,2,27,178.9620202823623),
    @(44531,4,21,139.1926824418373),
    @(44532,18,35,231.9878040697289),
    @(44533,8,40,265.128918936833),
    @(44534,6,45,298.2700338039372),
    @(44535,5,44,291.6418108305164),
    @(44536,3,46,304.898256777358),
    @(44537,6,50,331.4111486710413),
    @(44538,3,49,324.7829256976205)
)

$startRow = 386
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $vals = $newData[$i]
    $ws.Cells.Item($r, 1).Value2 = $vals[0]
    $ws.Cells.Item($r, 2).Value2 = $vals[1]
    $ws.Cells.Item($r, 3).Value2 = $vals[2]
    $ws.Cells.Item($r, 4).Value2 = $vals[3]
}

# Column A carries the same date style ("YYYY-MM-DD HH:MM:SS", centered, bordered)
# as the rest of the series above -- copy it down from the last pre-existing row.
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)

Write-Output "Added rows 386-464"